$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9289177060127258
$ws.Range("B1").Value = 1.67980694770813
$ws.Range("D1").Value = 1.840297698974609
$ws.Range("E1").Value = 1.090541839599609
